# PPT => NG 8
#
# 1) Remove the slide "Chargement des scripts" (the SystemJS script-loading
#    slide is no longer relevant once the deck moves on to Angular 8 /
#    Angular CLI, which doesn't use SystemJS).
# 2) On the "Ajout d'une expression" slide, the paragraph that used to be
#    typed as three separate runs ("{{}} " / "signifie " / "une expression")
#    is retyped as a single run with the same text.

$p = $ppt.ActivePresentation

# --- 1) Delete the "Chargement des scripts" slide -------------------------
$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
    if ($title -eq "Chargement des scripts") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -ge 1) {
    $p.Slides.Item($targetIndex).Delete()
}

# --- 2) Merge the "{{}} signifie une expression" runs into one run --------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            $full = $tr.Text
            if ($full -like "*signifie*") {
                $needle = "{{}} signifie une expression"
                $idx = $full.IndexOf($needle)
                if ($idx -ge 0) {
                    $sub = $tr.Characters($idx + 1, $needle.Length)
                    $sub.Text = $needle
                }
            }
        }
    }
}
